$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Client")

# Update the three client name cells (shared strings 0107A7/A8/A9 -> 0407A1/A2/A3)
$ws.Range("B2").Value = "Anh Tester Client 0407A1"
$ws.Range("B3").Value = "Anh Tester Client 0407A2"
$ws.Range("B4").Value = "Anh Tester Client 0407A3"

# Move the active selection on the Client sheet from B6 to B7
$ws.Activate()
$ws.Range("B7").Select()
